# Insert a new "School" column before the existing "Programs" column (old column G),
# shifting "Programs" to column H. Populate the new column with the school each
# child's family/guardian is associated with, then fix up column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G; this pushes the old G (Programs) to H and shifts
# all the cell data/formatting along with it.
$ws.Columns.Item(7).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 7).Value = "School"

# Populate the School column for each data row.
for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 7).Value = "Penn Hills Middle School"
}
for ($row = 24; $row -le 54; $row++) {
    $ws.Cells.Item($row, 7).Value = "Sunnyside"
}

# Column widths: School gets a wider column; Programs (now H) is left alone so it
# keeps the exact width it had as the original column G.
$ws.Columns.Item(7).ColumnWidth = 20
